# Commit work that I forgot to commit
#
# For the "Think of story" / "Draw basic textures" / "Write a list of
# required textures" / "Draw" / "Code prototype game" / "Design one level"
# list items:
#   - drop the explicit <w:jc w:val="left"/> from the paragraph properties
#   - (all but "Code prototype game") clear the run-level <w:rPr>, which only
#     held a redundant <w:shd/> already present on the paragraph mark's
#     <w:rPr>, down to an empty <w:rPr/>

$d = $word.ActiveDocument

# Paragraphs whose run carries a redundant shading <w:rPr> that needs to
# collapse to <w:rPr/>. "Code prototype game" is deliberately excluded: its
# run already has an empty <w:rPr/>, and the diff only drops <w:jc/> there.
$shadingTargets = @(
    "Think of story",
    "Draw basic textures",
    "Write a list of required textures",
    "Draw",
    "Design one level"
)

# All six paragraphs lose their explicit <w:jc w:val="left"/>.
$alignmentTargets = $shadingTargets + @("Code prototype game")

function Get-ParagraphByText($doc, $text) {
    foreach ($para in $doc.Paragraphs) {
        $t = $para.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $para
        }
    }
    return $null
}

foreach ($target in $shadingTargets) {
    $p = Get-ParagraphByText $d $target
    if ($null -eq $p) {
        throw "paragraph not found: '$target'"
    }

    $openxml = $p.Range.WordOpenXML
    if ($openxml -notmatch '(?s)<w:body>(<w:p[ >].*?</w:p>)') {
        throw "could not isolate paragraph XML for '$target'"
    }
    $paraXml = $matches[1]

    # The mini-document round trip stamps a w14:paraId/w:rsid* cluster onto
    # the <w:p> start tag; the source document never had those, so drop them
    # again before writing the paragraph back.
    $paraXml = $paraXml -replace '^<w:p\s+[^>]*>', '<w:p>'

    $marker = "</w:pPr>"
    $splitIdx = $paraXml.IndexOf($marker)
    if ($splitIdx -lt 0) {
        throw "no <w:pPr> found for '$target'"
    }
    $head = $paraXml.Substring(0, $splitIdx + $marker.Length)
    $tail = $paraXml.Substring($splitIdx + $marker.Length)

    # Remove the explicit left-justify from the paragraph properties.
    $head = $head -replace '<w:jc[^>]*/>', ''

    # Clear the run's own rPr (first rPr after </w:pPr>) to <w:rPr/>.
    $tail = $tail -replace '<w:rPr>.*?</w:rPr>', '<w:rPr/>'

    $newXml = $head + $tail
    $p.Range.InsertXML($newXml)
}

# "Code prototype game" (and, redundantly but harmlessly, the five
# paragraphs above since InsertXML already dropped their <w:jc/>): remove
# the explicit left alignment via the paragraph-format object model, which
# edits the <w:pPr> in place without disturbing anything else.
foreach ($target in $alignmentTargets) {
    $p = Get-ParagraphByText $d $target
    if ($null -eq $p) {
        throw "paragraph not found: '$target'"
    }
    $p.Format.Alignment = 0
}
